# Auto-generated Excel COM-interop script to apply numeric value updates
# to the Shiva_Profits workbook's 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 10705.667
$ws.Range("J40").Value = 7411.25
$ws.Range("L40").Value = 7411.25
$ws.Range("N40").Value = -7761.25
$ws.Range("H54").Value = 14474.75
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H69").Value = 32995.668
$ws.Range("I69").Value = 32996.5
$ws.Range("J69").Value = 32994
$ws.Range("K69").Value = 98989.5
$ws.Range("L69").Value = 98982
$ws.Range("M69").Value = -98115.5
$ws.Range("N69").Value = -100730
$ws.Range("H72").Value = 32995.668
$ws.Range("I72").Value = 32996.5
$ws.Range("J72").Value = 32994
$ws.Range("K72").Value = 296968.5
$ws.Range("L72").Value = 296946
$ws.Range("M72").Value = -292600.5
$ws.Range("N72").Value = -305682
$ws.Range("H82").Value = 1600
$ws.Range("I82").Value = 1600
$ws.Range("K82").Value = 4800
$ws.Range("M82").Value = -4394
$ws.Range("H85").Value = 1600
$ws.Range("I85").Value = 1600
$ws.Range("K85").Value = 4800
$ws.Range("M85").Value = -3396
$ws.Range("H88").Value = 8259.833000000001
$ws.Range("I88").Value = 3662.75
$ws.Range("J88").Value = 12856.917
$ws.Range("K88").Value = 3662.75
$ws.Range("L88").Value = 12856.917
$ws.Range("M88").Value = -3256.75
$ws.Range("N88").Value = -13668.917
$ws.Range("H91").Value = 8259.833000000001
$ws.Range("I91").Value = 3662.75
$ws.Range("J91").Value = 12856.917
$ws.Range("K91").Value = 3662.75
$ws.Range("L91").Value = 12856.917
$ws.Range("M91").Value = -2258.75
$ws.Range("N91").Value = -15664.917
$ws.Range("H112").Value = 2756.7632
$ws.Range("J112").Value = 3039.276
$ws.Range("L112").Value = 9117.828
$ws.Range("N112").Value = -11333.828
$ws.Range("H113").Value = 2552.8
$ws.Range("J113").Value = 2720.5
$ws.Range("L113").Value = 2720.5
$ws.Range("N113").Value = -9228.5
$ws.Range("H121").Value = 2820
$ws.Range("J121").Value = 2821.6875
$ws.Range("L121").Value = 8465.0625
$ws.Range("N121").Value = -11959.0625
$ws.Range("H132").Value = 24107.49
$ws.Range("I132").Value = 4654.108
$ws.Range("K132").Value = 13962.324
$ws.Range("M132").Value = -11432.324
$ws.Range("H135").Value = 1701.5652
$ws.Range("I135").Value = 1474.75
$ws.Range("K135").Value = 13272.75
$ws.Range("M135").Value = -10737.75
$ws.Range("H137").Value = 2468.5085
$ws.Range("I137").Value = 2433.4849
$ws.Range("K137").Value = 7300.4547
$ws.Range("M137").Value = -4750.4547
$ws.Range("H138").Value = 3695.5615
$ws.Range("I138").Value = 1898.4166
$ws.Range("J138").Value = 4575.796
$ws.Range("K138").Value = 5695.2498
$ws.Range("L138").Value = 13727.388
$ws.Range("M138").Value = -555.2497999999996
$ws.Range("N138").Value = -24007.388

# ---- Sheet 2: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1321.7903
$ws.Range("I2").Value = 1277.9818
$ws.Range("J2").Value = 1666
$ws.Range("K2").Value = 1277.9818
$ws.Range("L2").Value = 1666
$ws.Range("M2").Value = -1164.9818
$ws.Range("N2").Value = -1892
$ws.Range("H4").Value = 585.1429000000001
$ws.Range("I4").Value = 553.2308
$ws.Range("K4").Value = 553.2308
$ws.Range("M4").Value = -437.2308
$ws.Range("H5").Value = 362.5
$ws.Range("I5").Value = 362.5
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 362.5
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -250.5
$ws.Range("N5").ClearContents()
$ws.Range("H10").Value = 7698.3335
$ws.Range("I10").Value = 100
$ws.Range("J10").Value = 11497.5
$ws.Range("K10").Value = 100
$ws.Range("L10").Value = 11497.5
$ws.Range("M10").Value = 70
$ws.Range("N10").Value = -11837.5
$ws.Range("H28").Value = 21048.646
$ws.Range("I28").Value = 8261
$ws.Range("K28").Value = 8261
$ws.Range("M28").Value = -8069
$ws.Range("H32").Value = 2009.5172
$ws.Range("I32").Value = 1592.2927
$ws.Range("K32").Value = 1592.2927
$ws.Range("M32").Value = -1305.2927
$ws.Range("H34").Value = 40519.2
$ws.Range("I34").Value = 40519.2
$ws.Range("K34").Value = 40519.2
$ws.Range("M34").Value = -40248.2
$ws.Range("H99").Value = 21048.646
$ws.Range("I99").Value = 8261
$ws.Range("K99").Value = 8261
$ws.Range("M99").Value = -5266
$ws.Range("H110").Value = 2191.6
$ws.Range("I110").Value = 2047.9166
$ws.Range("K110").Value = 2047.9166
$ws.Range("M110").Value = -2.916600000000017
$ws.Range("H116").Value = 1321.7903
$ws.Range("I116").Value = 1277.9818
$ws.Range("J116").Value = 1666
$ws.Range("K116").Value = 1277.9818
$ws.Range("L116").Value = 1666
$ws.Range("M116").Value = 1016.0182
$ws.Range("N116").Value = -6254
$ws.Range("H122").Value = 6405.3516
$ws.Range("I122").Value = 4582.5
$ws.Range("K122").Value = 13747.5
$ws.Range("M122").Value = -11297.5
$ws.Range("H132").Value = 5719.679
$ws.Range("I132").Value = 2686.875
$ws.Range("K132").Value = 8060.625
$ws.Range("M132").Value = -5530.625
$ws.Range("H138").Value = 234947.56
$ws.Range("J138").Value = 234947.56
$ws.Range("L138").Value = 234947.56
$ws.Range("N138").Value = -245227.56

# ---- Sheet 3: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1321.7903
$ws.Range("I3").Value = 1277.9818
$ws.Range("J3").Value = 1666
$ws.Range("K3").Value = 1277.9818
$ws.Range("L3").Value = 1666
$ws.Range("M3").Value = -1163.9818
$ws.Range("N3").Value = -1894
$ws.Range("H4").Value = 362.5
$ws.Range("I4").Value = 362.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 362.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -247.5
$ws.Range("N4").ClearContents()
$ws.Range("H9").Value = 31000
$ws.Range("J9").Value = 31000
$ws.Range("L9").Value = 31000
$ws.Range("N9").Value = -31336
$ws.Range("H22").Value = 1418.75
$ws.Range("I22").Value = 1418.75
$ws.Range("K22").Value = 1418.75
$ws.Range("M22").Value = -1245.75
$ws.Range("H80").Value = 1852428.5
$ws.Range("J80").Value = 2564617.5
$ws.Range("L80").Value = 2564617.5
$ws.Range("N80").Value = -2566613.5
$ws.Range("H83").Value = 1852428.5
$ws.Range("J83").Value = 2564617.5
$ws.Range("L83").Value = 12823087.5
$ws.Range("N83").Value = -12833071.5
$ws.Range("H94").Value = 4339.5
$ws.Range("I94").Value = 3035.05
$ws.Range("J94").Value = 6948.4
$ws.Range("K94").Value = 3035.05
$ws.Range("L94").Value = 6948.4
$ws.Range("M94").Value = -2584.05
$ws.Range("N94").Value = -7850.4
$ws.Range("H97").Value = 13831.75
$ws.Range("I97").Value = 13831.75
$ws.Range("K97").Value = 13831.75
$ws.Range("M97").Value = -12840.75
$ws.Range("H99").Value = 5719.375
$ws.Range("I99").Value = 5876.1665
$ws.Range("J99").Value = 5249
$ws.Range("K99").Value = 5876.1665
$ws.Range("L99").Value = 5249
$ws.Range("M99").Value = -4378.1665
$ws.Range("N99").Value = -8245
$ws.Range("H105").Value = 1990.96
$ws.Range("I105").Value = 1926.4894
$ws.Range("K105").Value = 1926.4894
$ws.Range("M105").Value = -179.4893999999999
$ws.Range("H123").Value = 88999.75
$ws.Range("J123").Value = 88999.75
$ws.Range("L123").Value = 88999.75
$ws.Range("N123").Value = -98799.75
$ws.Range("H134").Value = 3578.4856
$ws.Range("I134").Value = 2835.3447
$ws.Range("K134").Value = 8506.034100000001
$ws.Range("M134").Value = -5971.034100000001

# ---- Sheet 4: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1335.2
$ws.Range("J16").Value = 1348.1923
$ws.Range("L16").Value = 1348.1923
$ws.Range("N16").Value = -1922.1923
$ws.Range("H68").Value = 99999
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 99999
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H94").Value = 4290.5625
$ws.Range("I94").Value = 3837.375
$ws.Range("K94").Value = 3837.375
$ws.Range("M94").Value = -3386.375
$ws.Range("H96").Value = 24811.334
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 24811.334
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 24811.334
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -30303.334
$ws.Range("H99").Value = 4370.5557
$ws.Range("J99").Value = 4379.5
$ws.Range("L99").Value = 4379.5
$ws.Range("N99").Value = -7375.5
$ws.Range("H113").Value = 1335.2
$ws.Range("J113").Value = 1348.1923
$ws.Range("L113").Value = 1348.1923
$ws.Range("N113").Value = -5688.1923
$ws.Range("H126").Value = 4370.5557
$ws.Range("J126").Value = 4379.5
$ws.Range("L126").Value = 13138.5
$ws.Range("N126").Value = -18078.5
$ws.Range("H134").Value = 2799.8845
$ws.Range("I134").Value = 1986.1628
$ws.Range("J134").Value = 6687.6665
$ws.Range("K134").Value = 5958.4884
$ws.Range("L134").Value = 20062.9995
$ws.Range("M134").Value = -3423.4884
$ws.Range("N134").Value = -25132.9995

# ---- Sheet 5: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9337.682000000001
$ws.Range("I3").Value = 3028.6
$ws.Range("K3").Value = 9085.799999999999
$ws.Range("M3").Value = -8973.799999999999
$ws.Range("H23").Value = 937.381
$ws.Range("J23").Value = 264.22223
$ws.Range("L23").Value = 792.66669
$ws.Range("N23").Value = -1262.66669
$ws.Range("H37").Value = 76614.57000000001
$ws.Range("J37").Value = 76614.57000000001
$ws.Range("L37").Value = 229843.71
$ws.Range("N37").Value = -230067.71
$ws.Range("H55").Value = 630.1818
$ws.Range("I55").Value = 116.5
$ws.Range("K55").Value = 349.5
$ws.Range("M55").Value = -172.5
$ws.Range("H68").Value = 1742.3572
$ws.Range("I68").Value = 771.1429000000001
$ws.Range("K68").Value = 2313.4287
$ws.Range("M68").Value = -1502.4287
$ws.Range("H71").Value = 1742.3572
$ws.Range("I71").Value = 771.1429000000001
$ws.Range("K71").Value = 6940.2861
$ws.Range("M71").Value = -2884.2861
$ws.Range("H106").Value = 10000
$ws.Range("J106").Value = 10000
$ws.Range("L106").Value = 30000
$ws.Range("N106").Value = -31892
$ws.Range("H116").Value = 2861.2727
$ws.Range("I116").Value = 3232.5715
$ws.Range("K116").Value = 9697.7145
$ws.Range("M116").Value = -6255.7145
$ws.Range("H136").Value = 3831.7334
$ws.Range("I136").Value = 2113.6155
$ws.Range("K136").Value = 6340.8465
$ws.Range("M136").Value = -1240.8465
$ws.Range("H137").Value = 3850.375
$ws.Range("I137").Value = 3407.7273
$ws.Range("J137").Value = 4824.2
$ws.Range("K137").Value = 10223.1819
$ws.Range("L137").Value = 14472.6
$ws.Range("M137").Value = -5123.1819
$ws.Range("N137").Value = -24672.6

# ---- Sheet 6: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3130823.8
$ws.Range("I11").Value = 3338001
$ws.Range("J11").Value = 3106918.8
$ws.Range("K11").Value = 3338001
$ws.Range("L11").Value = 3106918.8
$ws.Range("M11").Value = -3337862
$ws.Range("N11").Value = -3107196.8
$ws.Range("H18").Value = 49996.285
$ws.Range("J18").Value = 63324.668
$ws.Range("L18").Value = 63324.668
$ws.Range("N18").Value = -63910.668
$ws.Range("H102").Value = 5108.7954
$ws.Range("J102").Value = 3689.3635
$ws.Range("L102").Value = 3689.3635
$ws.Range("N102").Value = -6933.363499999999
$ws.Range("H126").Value = 5692.4707
$ws.Range("I126").Value = 5175.909
$ws.Range("K126").Value = 15527.727
$ws.Range("M126").Value = -13057.727
$ws.Range("H132").Value = 9021.909
$ws.Range("I132").Value = 8736.385
$ws.Range("K132").Value = 26209.155
$ws.Range("M132").Value = -23679.155
$ws.Range("H135").Value = 149131.61
$ws.Range("J135").Value = 149131.61
$ws.Range("L135").Value = 149131.61
$ws.Range("N135").Value = -159271.61

# ---- Sheet 7: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9479.700000000001
$ws.Range("I7").Value = 9185.143
$ws.Range("J7").Value = 10167
$ws.Range("K7").Value = 9185.143
$ws.Range("L7").Value = 10167
$ws.Range("M7").Value = -9073.143
$ws.Range("N7").Value = -10391
$ws.Range("H16").Value = 1367.8857
$ws.Range("I16").Value = 1446.3
$ws.Range("J16").Value = 897.4
$ws.Range("K16").Value = 1446.3
$ws.Range("L16").Value = 897.4
$ws.Range("M16").Value = -1276.3
$ws.Range("N16").Value = -1237.4
$ws.Range("H25").Value = 21054.428
$ws.Range("I25").Value = 17485.8
$ws.Range("K25").Value = 17485.8
$ws.Range("M25").Value = -17255.8
$ws.Range("H46").Value = 1617.625
$ws.Range("I46").Value = 1042.5454
$ws.Range("K46").Value = 1042.5454
$ws.Range("M46").Value = -854.5454
$ws.Range("H61").Value = 13428.954
$ws.Range("I61").Value = 14042.35
$ws.Range("K61").Value = 14042.35
$ws.Range("M61").Value = -13840.35
$ws.Range("H64").Value = 29250
$ws.Range("J64").Value = 23000
$ws.Range("L64").Value = 23000
$ws.Range("N64").Value = -23450
$ws.Range("H67").Value = 29250
$ws.Range("J67").Value = 23000
$ws.Range("L67").Value = 23000
$ws.Range("N67").Value = -24560
$ws.Range("H68").Value = 3185.913
$ws.Range("I68").Value = 2986.9412
$ws.Range("K68").Value = 2986.9412
$ws.Range("M68").Value = -2237.9412
$ws.Range("H70").Value = 29999
$ws.Range("J70").Value = 29999
$ws.Range("L70").Value = 29999
$ws.Range("N70").Value = -30539
$ws.Range("H71").Value = 3185.913
$ws.Range("I71").Value = 2986.9412
$ws.Range("K71").Value = 14934.706
$ws.Range("M71").Value = -11190.706
$ws.Range("H73").Value = 29999
$ws.Range("J73").Value = 29999
$ws.Range("L73").Value = 29999
$ws.Range("N73").Value = -31871
$ws.Range("H93").Value = 4518.154
$ws.Range("I93").Value = 4814.189
$ws.Range("J93").Value = 3787.9333
$ws.Range("K93").Value = 4814.189
$ws.Range("L93").Value = 3787.9333
$ws.Range("M93").Value = -3566.189
$ws.Range("N93").Value = -6283.933300000001
$ws.Range("H113").Value = 13428.954
$ws.Range("I113").Value = 14042.35
$ws.Range("K113").Value = 14042.35
$ws.Range("M113").Value = -11872.35
$ws.Range("H126").Value = 9479.700000000001
$ws.Range("I126").Value = 9185.143
$ws.Range("J126").Value = 10167
$ws.Range("K126").Value = 27555.429
$ws.Range("L126").Value = 30501
$ws.Range("M126").Value = -25085.429
$ws.Range("N126").Value = -35441
$ws.Range("H132").Value = 44741
$ws.Range("I132").Value = 145172.42
$ws.Range("J132").Value = 9590
$ws.Range("K132").Value = 435517.26
$ws.Range("L132").Value = 28770
$ws.Range("M132").Value = -432987.26
$ws.Range("N132").Value = -33830
$ws.Range("H136").Value = 7901.8945
$ws.Range("I136").Value = 2759.8572
$ws.Range("J136").Value = 22299.6
$ws.Range("K136").Value = 8279.571599999999
$ws.Range("L136").Value = 66898.79999999999
$ws.Range("M136").Value = -5729.571599999999
$ws.Range("N136").Value = -71998.79999999999
$ws.Range("H141").Value = 161623.38
$ws.Range("J141").Value = 161623.38
$ws.Range("L141").Value = 161623.38
$ws.Range("N141").Value = -171983.38

# ---- Sheet 8: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 3146803
$ws.Range("I4").Value = 4231304
$ws.Range("J4").Value = 13800.223
$ws.Range("K4").Value = 4231304
$ws.Range("L4").Value = 13800.223
$ws.Range("M4").Value = -4231191
$ws.Range("N4").Value = -14026.223
$ws.Range("H38").Value = 31000
$ws.Range("J38").Value = 31000
$ws.Range("L38").Value = 31000
$ws.Range("N38").Value = -31946
$ws.Range("H46").Value = 55201.883
$ws.Range("J46").Value = 55201.883
$ws.Range("L46").Value = 55201.883
$ws.Range("N46").Value = -55663.883
$ws.Range("H48").Value = 33000
$ws.Range("J48").Value = 33000
$ws.Range("L48").Value = 33000
$ws.Range("N48").Value = -34138
$ws.Range("H49").Value = 47499.5
$ws.Range("I49").Value = 47499.5
$ws.Range("K49").Value = 47499.5
$ws.Range("M49").Value = -47269.5
$ws.Range("H100").Value = 1115.2258
$ws.Range("I100").Value = 593.35
$ws.Range("K100").Value = 1186.7
$ws.Range("M100").Value = -645.7
$ws.Range("H107").Value = 1859.0625
$ws.Range("I107").Value = 1590.8
$ws.Range("K107").Value = 4772.4
$ws.Range("M107").Value = -2852.4
$ws.Range("H113").Value = 1818.125
$ws.Range("I113").Value = 1393.9412
$ws.Range("K113").Value = 4181.8236
$ws.Range("M113").Value = -2011.8236
$ws.Range("H122").Value = 6770.484
$ws.Range("I122").Value = 2314.1365
$ws.Range("J122").Value = 17663.777
$ws.Range("K122").Value = 6942.4095
$ws.Range("L122").Value = 52991.33099999999
$ws.Range("M122").Value = -4492.4095
$ws.Range("N122").Value = -57891.33099999999
$ws.Range("H126").Value = 3875.8215
$ws.Range("I126").Value = 3492.625
$ws.Range("J126").Value = 6175
$ws.Range("K126").Value = 10477.875
$ws.Range("L126").Value = 18525
$ws.Range("M126").Value = -8007.875
$ws.Range("N126").Value = -23465
$ws.Range("H132").Value = 13162.462
$ws.Range("I132").Value = 12611.5
$ws.Range("K132").Value = 37834.5
$ws.Range("M132").Value = -35304.5
$ws.Range("H134").Value = 55201.883
$ws.Range("J134").Value = 55201.883
$ws.Range("L134").Value = 165605.649
$ws.Range("N134").Value = -170675.649
$ws.Range("H136").Value = 5125.5757
$ws.Range("I136").Value = 6989.8945
$ws.Range("J136").Value = 2595.4285
$ws.Range("K136").Value = 20969.6835
$ws.Range("L136").Value = 7786.2855
$ws.Range("M136").Value = -18419.6835
$ws.Range("N136").Value = -12886.2855
